$d = $word.ActiveDocument

# --- Change 1 ------------------------------------------------------------
# " In verschillende assignments hebben we voorspellingen gemaakt, ..."
# becomes
# " In verschillende assignments (zoals assignment 19 met clusters) hebben
#   we voorspellingen gemaakt, ..."
$d.Content.Find.Execute(
    "In verschillende assignments hebben we voorspellingen",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In verschillende assignments (zoals assignment 19 met clusters) hebben we voorspellingen",
    2) | Out-Null

# --- Change 2 --------------------------------------------------------------
# Insert three brand new paragraphs right after the "Data scientists
# rapporteren en delen verder ook de gegenereerde inzichten..." paragraph,
# and before the two trailing blank paragraphs at the end of the document.

# That paragraph (ending in "...of andere betrokkenen.") is paragraph 13.
$lastPara = $d.Paragraphs(13).Range
$lastPara.InsertParagraphAfter()

$p1 = $d.Paragraphs(14).Range
$p1.InsertBefore("Bij de eerste vier jupyter notebook assignments heb je een voorbeeld van informatie extraheren uit data (univariate analysis). We wilden bijvoorbeeld zien wat de verdeelsleutel is van de geslachten bij de pinguïns. Met barplots, pie charts en linecharts kunnen we die verdeelsleutel heel overzichtelijk zien.")
$d.Paragraphs(14).Range.InsertParagraphAfter()

$p2 = $d.Paragraphs(15).Range
$p2.InsertBefore("Dit geld hetzelfde voor de bivariate analysis. Doordat je meerdere kolommen hebt, zie je eerder een correlatie tussen de kolommen. In assignment 9 en 10 heb je bijvoorbeeld negatieve- en positieve correlaties. Hiermee kan je mogelijk achter de reden van de correlatie komen.")
$d.Paragraphs(15).Range.InsertParagraphAfter()

$p3 = $d.Paragraphs(16).Range
$p3.InsertBefore("Als laatst kan een data scientist een besluit maken of tenminste advies geven op basis van de voorspellingen of ontdekte inzichten.")
